# Add two new columns "I0" (I) and "IF" (J) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
# Copy the formatting from the existing header cell H1 (bold, bordered,
# centered style) onto the two new header cells so they reuse the same
# cell style as the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-39) : [row, I value, J value] ---
$data = @(
    @(2, 9, 9),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 2, 4),
    @(8, 6, 7),
    @(9, 2, 5),
    @(10, 9, 9),
    @(11, 6, 8),
    @(12, 7, 8),
    @(13, 8, 8),
    @(14, 8, 9),
    @(15, 10, 10),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 6, 7),
    @(19, 6, 7),
    @(20, 7, 7),
    @(21, 6, 7),
    @(22, 7, 7),
    @(23, 5, 6),
    @(24, 6, 6),
    @(25, 6, 7),
    @(26, 7, 7),
    @(27, 6, 6),
    @(28, 7, 7),
    @(29, 6, 7),
    @(30, 8, 9),
    @(31, 8, 8),
    @(32, 6, 7),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 4, 5),
    @(36, 5, 6),
    @(37, 4, 4),
    @(38, 5, 5),
    @(39, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Host "Added I0/IF columns"
